# Generate Report for Handback
# The f4932e3c-...md file has now been handed back (in sync with en-US),
# so update its status on every sheet and record the new handback
# timestamps for the zh-cn and de-de locales.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet: row 3 is the f4932e3c-...md file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack

# --- zh-cn sheet: row 3 is the f4932e3c-...md file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusHandedBack
$wsZhCn.Range("H3").Value = "2016-03-21 08:42:42"

# --- de-de sheet: row 3 is the f4932e3c-...md file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusHandedBack
$wsDeDe.Range("H3").Value = "2016-03-21 08:42:48"
